$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns remain text, matching the
# original inline-string cell type (values include formats like "1.00",
# "96.784.50" and percentage strings that Excel would otherwise reinterpret).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "97.222.19"
$ws.Range("E2").Value = "  -1.65%  "

$ws.Range("D3").Value = "3.331.96"
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "246.26"
$ws.Range("E5").Value = "  -4.94%  "

$ws.Range("D6").Value = "649.78"
$ws.Range("E6").Value = "  -2.87%  "

$ws.Range("D7").Value = "1.35"
$ws.Range("E7").Value = "  -12.62%  "

$ws.Range("D8").Value = "0.416"
$ws.Range("E8").Value = "  -11.37%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "0.993"
$ws.Range("E10").Value = "  -8.68%  "

$ws.Range("D11").Value = "3.330.96"
$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("D12").Value = "0.204"
$ws.Range("E12").Value = "  -6.83%  "

$ws.Range("D13").Value = "39.86"
$ws.Range("E13").Value = "  -7.27%  "

$ws.Range("D14").Value = "96.783.96"
$ws.Range("E14").Value = "  -2.17%  "

$ws.Range("D15").Value = "5.99"
$ws.Range("E15").Value = "  +2.89%  "

$ws.Range("D16").Value = "0.0000251"
$ws.Range("E16").Value = "  -9.28%  "

$ws.Range("D17").Value = "3.950.72"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "8.39"
$ws.Range("E18").Value = "  +4.11%  "

$ws.Range("D19").Value = "3.324.52"
$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("D20").Value = "0.537"
$ws.Range("E20").Value = "  +22.17%  "

$ws.Range("D21").Value = "16.65"
$ws.Range("E21").Value = "  -3.91%  "

$ws.Range("D22").Value = "10.52"
$ws.Range("E22").Value = "  -0.77%  "

$ws.Range("D23").Value = "492.63"
$ws.Range("E23").Value = "  -7.73%  "

$ws.Range("D24").Value = "3.28"
$ws.Range("E24").Value = "  -7.83%  "

$ws.Range("D25").Value = "0.0000197"
$ws.Range("E25").Value = "  -9.88%  "

$ws.Range("D26").Value = "6.28"
$ws.Range("E26").Value = "  -1.74%  "

$ws.Range("D27").Value = "92.41"
$ws.Range("E27").Value = "  -10.18%  "

$ws.Range("D28").Value = "11.94"
$ws.Range("E28").Value = "  -6.61%  "

$ws.Range("D29").Value = "3.495.32"
$ws.Range("E29").Value = "  -2.78%  "

$ws.Range("D30").Value = "0.144"
$ws.Range("E30").Value = "  -4.86%  "

$ws.Range("D31").Value = "0.994"
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("D32").Value = "10.83"
$ws.Range("E32").Value = "  -6.96%  "

$ws.Range("D33").Value = "0.188"
$ws.Range("E33").Value = "  -4.03%  "

$ws.Range("D34").Value = "2.44"
$ws.Range("E34").Value = "  +11.89%  "

$ws.Range("D35").Value = "0.992"
$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("D36").Value = "0.542"
$ws.Range("E36").Value = "  -3.13%  "

$ws.Range("D37").Value = "28.11"
$ws.Range("E37").Value = "  -7.93%  "

$ws.Range("D38").Value = "7.53"
$ws.Range("E38").Value = "  -3.65%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "1.43"
$ws.Range("E39").Value = "  +4.91%  "

$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").Value = "0.149"
$ws.Range("E41").Value = "  -8.82%  "

$ws.Range("D42").Value = "498.53"
$ws.Range("E42").Value = "  -5.59%  "

$ws.Range("D43").Value = "24.54"
$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("D44").Value = "3.71"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("D45").Value = "0.825"
$ws.Range("E45").Value = "  -3.11%  "

$ws.Range("D46").Value = "8.56"
$ws.Range("E46").Value = "  +5.08%  "

$ws.Range("D47").Value = "0.0405"
$ws.Range("E47").Value = "  -6.40%  "

$ws.Range("D48").Value = "5.44"
$ws.Range("E48").Value = "  +4.75%  "

$ws.Range("D49").Value = "1.62"
$ws.Range("E49").Value = "  +3.65%  "

$ws.Range("D50").Value = "52.81"
$ws.Range("E50").Value = "  +4.05%  "

$ws.Range("D51").Value = "3.12"
$ws.Range("E51").Value = "  -11.63%  "
